$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the header row (row 1), shifting the
# existing data rows (2-22) down to rows 4-24.
$ws.Rows.Item(2).Resize(2).Insert()

# The inserted rows pick up the header row's formatting (bold/border); the
# data rows in this sheet carry no explicit style, so strip it back off.
$ws.Rows.Item(2).Resize(2).ClearFormats()

# Populate the two newly inserted rows with the new sensor readings.
$ws.Range("A2").Value = 0.0218384321779012
$ws.Range("B2").Value = -0.036499198526144
$ws.Range("C2").Value = -0.0226020142436027

$ws.Range("A3").Value = -0.00137444678694
$ws.Range("B3").Value = 0.0125227374956011
$ws.Range("C3").Value = 0.0018325957935303

# The last three rows of the original data (now at rows 22-24 after the
# insert above) are no longer part of the dataset, so remove them.
$ws.Rows.Item(22).Resize(3).Delete()
